# Make username and password in commands file generic
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace the hard-coded username value "gracz" with a generic placeholder
$ws.Range("C6").Value = "<username>"

# Replace the hard-coded password value "mattsNewPassword!" with a generic placeholder
$ws.Range("C11").Value = "<password>"

# Update the active selection to C11 (last edited cell)
$ws.Range("C11").Select()
